# Weekly update: a new sampling row is inserted into the Mango price series
# at row 62, pushing the existing rows (62-94) down by one (to 63-95).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 62; Excel shifts everything
# below down by one (rows 62..94 become 63..95) and extends the used range.
$ws.Rows("62").Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(62, 1).Value = 5
$ws.Cells.Item(62, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(62, 3).Value = "Maule"
$ws.Cells.Item(62, 4).Value = 44529
$ws.Cells.Item(62, 5).Value = 7
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100108
$ws.Cells.Item(62, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(62, 9).Value = 100108002
$ws.Cells.Item(62, 10).Value = "Mango"
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 220
$ws.Cells.Item(62, 14).Value = 6000
$ws.Cells.Item(62, 15).Value = 6000
$ws.Cells.Item(62, 16).Value = 6000
$ws.Cells.Item(62, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(62, 18).Value = "Perú"
$ws.Cells.Item(62, 19).Value = 1500
$ws.Cells.Item(62, 20).Value = 4
